# Arrange UML folders and add sequence diagram paths
#
# The lone "Click to edit Master title style" Title placeholder on the
# enumerations/playerTurn diagram slide was removed (it carried no real
# content), and the whole UML diagram group was shifted up to take over
# the vertical space that placeholder used to occupy.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Precise target "Top" (in points) for every remaining shape, keyed by the
# shape's (stable) Id. These correspond to shifting every shape's EMU
# y-offset up by 1173892 EMU, computed precisely (EMU / 12700 + a half-EMU
# epsilon so float32 marshalling lands back on the exact integer EMU).
$targetTop = @{
  54 = 129.56759842519688
  53 = 129.56759842519688
  5  = 153.56759842519688
  6  = 247.026968503937
  7  = 153.56759842519688
  8  = 168.1082283464567
  9  = 214.8379133858268
  10 = 272.64885826771655
  11 = 160.6082283464567
  12 = 140.64885826771655
  13 = 218.64885826771655
  14 = 212.64885826771655
  15 = 231.56759842519688
  16 = 297.5675984251969
  17 = 339.5675984251969
  23 = 154.48633858267718
  24 = 244.48633858267718
  25 = 154.48633858267718
  26 = 169.026968503937
  27 = 212.526968503937
  28 = 273.5676771653543
  29 = 171.79728346456693
  30 = 141.56759842519688
  31 = 219.56759842519688
  32 = 213.56759842519688
  33 = 261.5675984251969
  49 = 323.58082677165356
  50 = 372.6753937007874
  55 = 261.5675984251969
}

# 1) Delete the empty title placeholder.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Type -eq 14 -and $shp.Name -eq "Title 37") {
        $shp.Delete()
    }
}

# 2) Slide the rest of the UML diagram up into the freed space.
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($targetTop.ContainsKey($shp.Id)) {
        $shp.Top = $targetTop[$shp.Id]
    }
}
